# error solve ifrs list
# Rewrites the numeric financial data cells (columns D:AJ) for data rows 2-6
# with corrected figures, and clears out the erroneous trailing rows 7-9
# (which had been populated with stale/incorrect duplicate data) down to
# just their identifying columns A:C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (2014/12) ---
$ws.Range("D2").Value = 2108
$ws.Range("E2").Value = 56
$ws.Range("F2").Value = 56
$ws.Range("G2").Value = 47
$ws.Range("H2").Value = 36
$ws.Range("I2").Value = 37
$ws.Range("J2").Value = -1
$ws.Range("K2").Value = 1984
$ws.Range("L2").Value = 765
$ws.Range("M2").Value = 1218
$ws.Range("N2").Value = 1215
$ws.Range("O2").Value = 4
$ws.Range("P2").Value = 125
$ws.Range("Q2").Value = 207
$ws.Range("R2").Value = -58
$ws.Range("S2").Value = -122
$ws.Range("T2").Value = 21
$ws.Range("U2").Value = 186
$ws.Range("V2").Value = 95
$ws.Range("W2").Value = 2.64
$ws.Range("X2").Value = 1.71
$ws.Range("Y2").Value = 3.01
$ws.Range("Z2").Value = 1.79
$ws.Range("AA2").Value = 62.79
$ws.Range("AB2").Value = 904.74
$ws.Range("AC2").Value = 146
$ws.Range("AD2").Value = 21.2
$ws.Range("AE2").Value = 5218
$ws.Range("AF2").Value = 0.59
$ws.Range("AG2").Value = 100
$ws.Range("AH2").Value = 3.23
$ws.Range("AI2").Value = 63.7
$ws.Range("AJ2").Value = 25000000

# --- Row 3 (2015/12) ---
$ws.Range("D3").Value = 2669
$ws.Range("E3").Value = 123
$ws.Range("F3").Value = 123
$ws.Range("G3").Value = 117
$ws.Range("H3").Value = 82
$ws.Range("I3").Value = 83
$ws.Range("J3").Value = -2
$ws.Range("K3").Value = 2347
$ws.Range("L3").Value = 1101
$ws.Range("M3").Value = 1246
$ws.Range("N3").Value = 1244
$ws.Range("O3").Value = 2
$ws.Range("P3").Value = 125
$ws.Range("Q3").Value = 452
$ws.Range("R3").Value = -240
$ws.Range("S3").Value = -69
$ws.Range("T3").Value = 11
$ws.Range("U3").Value = 441
$ws.Range("V3").Value = 76
$ws.Range("W3").Value = 4.61
$ws.Range("X3").Value = 3.07
$ws.Range("Y3").Value = 6.79
$ws.Range("Z3").Value = 3.78
$ws.Range("AA3").Value = 88.36
$ws.Range("AB3").Value = 949.66
$ws.Range("AC3").Value = 334
$ws.Range("AD3").Value = 10.71
$ws.Range("AE3").Value = 5529
$ws.Range("AF3").Value = 0.65
$ws.Range("AG3").Value = 100
$ws.Range("AH3").Value = 2.8
$ws.Range("AI3").Value = 26.95
$ws.Range("AJ3").Value = 25000000

# --- Row 4 (2016/12) ---
$ws.Range("D4").Value = 2745
$ws.Range("E4").Value = 222
$ws.Range("F4").Value = 222
$ws.Range("G4").Value = 202
$ws.Range("H4").Value = 160
$ws.Range("I4").Value = 161
$ws.Range("J4").Value = -1
$ws.Range("K4").Value = 2454
$ws.Range("L4").Value = 976
$ws.Range("M4").Value = 1478
$ws.Range("N4").Value = 1475
$ws.Range("O4").Value = 4
$ws.Range("P4").Value = 125
$ws.Range("Q4").Value = 168
$ws.Range("R4").Value = -348
$ws.Range("S4").Value = 72
$ws.Range("T4").Value = 18
$ws.Range("U4").Value = 150
$ws.Range("V4").Value = 75
$ws.Range("W4").Value = 8.08
$ws.Range("X4").Value = 5.83
$ws.Range("Y4").Value = 11.82
$ws.Range("Z4").Value = 6.67
$ws.Range("AA4").Value = 66.04000000000001
$ws.Range("AB4").Value = 1075.36
$ws.Range("AC4").Value = 642
$ws.Range("AD4").Value = 5.45
$ws.Range("AE4").Value = 5898
$ws.Range("AF4").Value = 0.59
$ws.Range("AG4").Value = 130
$ws.Range("AH4").Value = 3.71
$ws.Range("AI4").Value = 20.23
$ws.Range("AJ4").Value = 25000000

# --- Row 5 (2017/12) ---
$ws.Range("D5").Value = 2657
$ws.Range("E5").Value = 242
$ws.Range("F5").Value = 242
$ws.Range("G5").Value = 284
$ws.Range("H5").Value = 214
$ws.Range("I5").Value = 215
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 2564
$ws.Range("L5").Value = 887
$ws.Range("M5").Value = 1677
$ws.Range("N5").Value = 1645
$ws.Range("O5").Value = 32
$ws.Range("P5").Value = 125
$ws.Range("Q5").Value = 55
$ws.Range("R5").Value = -35
$ws.Range("S5").Value = -19
$ws.Range("T5").Value = 11
$ws.Range("U5").Value = 43
$ws.Range("V5").Value = 75
$ws.Range("W5").Value = 9.109999999999999
$ws.Range("X5").Value = 8.07
$ws.Range("Y5").Value = 13.77
$ws.Range("Z5").Value = 8.539999999999999
$ws.Range("AA5").Value = 52.89
$ws.Range("AB5").Value = 1222.56
$ws.Range("AC5").Value = 859
$ws.Range("AD5").Value = 5.32
$ws.Range("AE5").Value = 6687
$ws.Range("AF5").Value = 0.68
$ws.Range("AG5").Value = 160
$ws.Range("AH5").Value = 3.5
$ws.Range("AI5").Value = 18.38
$ws.Range("AJ5").Value = 25000000

# --- Row 6 (2018/12) --- (no J6/O6 in this sheet, stays absent)
$ws.Range("D6").Value = 2987
$ws.Range("E6").Value = 244
$ws.Range("F6").Value = 244
$ws.Range("G6").Value = 277
$ws.Range("H6").Value = 206
$ws.Range("I6").Value = 208
$ws.Range("K6").Value = 2769
$ws.Range("L6").Value = 928
$ws.Range("M6").Value = 1841
$ws.Range("N6").Value = 1797
$ws.Range("P6").Value = 125
$ws.Range("Q6").Value = 70
$ws.Range("R6").Value = -18
$ws.Range("S6").Value = -68
$ws.Range("T6").Value = 19
$ws.Range("U6").Value = 51
$ws.Range("V6").Value = 75
$ws.Range("W6").Value = 8.17
$ws.Range("X6").Value = 6.91
$ws.Range("Y6").Value = 12.09
$ws.Range("Z6").Value = 7.74
$ws.Range("AA6").Value = 50.43
$ws.Range("AB6").Value = 1379.38
$ws.Range("AC6").Value = 832
$ws.Range("AD6").Value = 6.06
$ws.Range("AE6").Value = 7489
$ws.Range("AF6").Value = 0.67
$ws.Range("AG6").Value = 160
$ws.Range("AH6").Value = 3.17
$ws.Range("AI6").Value = 18.49
$ws.Range("AJ6").Value = 25000000

# --- Rows 7-9 (2019/12(E), 2020/12(E), 2021/12(E)) ---
# These rows were erroneously populated with duplicate/incorrect figures;
# clear the data columns (D:AJ) back out, leaving only the A/B/C
# identifier columns intact.
$ws.Range("D7:AJ9").ClearContents()
